$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '51.729.96'
$ws.Range('E2').Value2 = '  +1.64%  '
$ws.Range('D3').Value2 = '2.798.61'
$ws.Range('E3').Value2 = '  +2.51%  '
$ws.Range('E4').Value2 = '  -0.31%  '
$ws.Range('D5').Value2 = '353.12'
$ws.Range('E5').Value2 = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '111.70'
$ws.Range('E6').Value2 = '  +5.29%  '
$ws.Range('E7').Value2 = '  +2.36%  '
$ws.Range('E8').Value2 = '  -0.15%  '
$ws.Range('D9').Value2 = '0.623'
$ws.Range('E9').Value2 = '  +9.14%  '
$ws.Range('D10').Value2 = '40.31'
$ws.Range('E10').Value2 = '  +4.46%  '
$ws.Range('E11').Value2 = '  -0.02%  '
$ws.Range('E12').Value2 = '  +1.58%  '
$ws.Range('E13').Value2 = '  +2.46%  '
$ws.Range('E14').Value2 = '  +5.12%  '
$ws.Range('D15').Value2 = '3.236.99'
$ws.Range('E15').Value2 = '  +1.49%  '
$ws.Range('D16').Value2 = '2.800.83'
$ws.Range('E16').Value2 = '  +1.67%  '
$ws.Range('D17').Value2 = '0.943'
$ws.Range('E17').Value2 = '  +3.70%  '
$ws.Range('D18').Value2 = '51.691.23'
$ws.Range('E18').Value2 = '  +1.37%  '
$ws.Range('D19').Value2 = '7.61'
$ws.Range('E19').Value2 = '  +1.38%  '
$ws.Range('D20').Value2 = '3.21'
$ws.Range('E20').Value2 = '  +7.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '13.60'
$ws.Range('E21').Value2 = '  +6.38%  '
$ws.Range('E22').Value2 = '  +2.68%  '
$ws.Range('D23').Value2 = '70.26'
$ws.Range('E23').Value2 = '  +2.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '267.30'
$ws.Range('E24').Value2 = '  +2.47%  '
$ws.Range('E25').Value2 = '  +2.52%  '
$ws.Range('E26').Value2 = '  -0.12%  '
$ws.Range('D27').Value2 = '26.13'
$ws.Range('E27').Value2 = '  +2.07%  '
$ws.Range('E28').Value2 = '  +0.26%  '
$ws.Range('D29').Value2 = '39.04'
$ws.Range('E29').Value2 = '  +15.14%  '
$ws.Range('E30').Value2 = '  +4.36%  '
$ws.Range('E31').Value2 = '  +0.84%  '
$ws.Range('D32').Value2 = '52.48'
$ws.Range('E32').Value2 = '  +2.06%  '
$ws.Range('D33').Value2 = '6.14'
$ws.Range('E33').Value2 = '  +3.34%  '
$ws.Range('D34').Value2 = '0.0453'
$ws.Range('E34').Value2 = '  +3.80%  '
$ws.Range('D35').Value2 = '0.0894'
$ws.Range('E35').Value2 = '  +8.94%  '
$ws.Range('E36').Value2 = '  +9.27%  '
$ws.Range('E37').Value2 = '  -0.43%  '
$ws.Range('D38').Value2 = '18.89'
$ws.Range('E38').Value2 = '  +4.46%  '
$ws.Range('D39').Value2 = '3.17'
$ws.Range('E39').Value2 = '  +2.19%  '
$ws.Range('E40').Value2 = '  +5.22%  '
$ws.Range('E41').Value2 = '  +2.55%  '
$ws.Range('D42').Value2 = '2.52'
$ws.Range('E42').Value2 = '  +3.82%  '
$ws.Range('E43').Value2 = '  +1.22%  '
$ws.Range('D44').Value2 = '120.21'
$ws.Range('E44').Value2 = '  +1.14%  '
$ws.Range('D45').Value2 = '21.92'
$ws.Range('E45').Value2 = '  +1.36%  '
$ws.Range('B46').Value2 = 'NEARProtocol'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value2 = '3.46'
$ws.Range('E46').Value2 = '  +9.46%  '
$ws.Range('B47').Value2 = 'ApeXProtocol'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value2 = '2.46'
$ws.Range('E47').Value2 = '  +6.47%  '
$ws.Range('D48').Value2 = '2.106.95'
$ws.Range('E48').Value2 = '  +2.36%  '
$ws.Range('D49').Value2 = '0.957'
$ws.Range('E49').Value2 = '  +6.39%  '
$ws.Range('D50').Value2 = '5.47'
$ws.Range('E50').Value2 = '  +1.13%  '
$ws.Range('E51').Value2 = '  +7.93%  '
